# Automatic update of files.
# Update the "Förändrad" (Modified) date column (C) for rows 2-42
# from serial date 45718 (2025-03-02) to 45719 (2025-03-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45718) {
        $cell.Value2 = 45719
    }
}
